# Update 13C-MFA files (run and result) for SC and IO under WT-batch and chemostats
#
# Changes applied:
#  1. FluxData: insert a new flux row "EX_glc__D_e.f" right after "BIOMASS.f"
#     (old row 3), shifting the remaining flux rows down by one; update the
#     (now shifted) old "EX_c5sugal_e.f" row's value/error, and append one
#     brand-new row "DIL_ade_d1.f" at the very end (row 38) with default
#     100 / 0.0001 placeholders; also fix BIOMASS.f's error from 1E-06 to
#     0.0001.
#  2. View-state refresh: zoom 60% -> 95% on every sheet, FluxData becomes
#     the active tab/sheet (was MSData), and selections collapse down to a
#     single cell on each sheet.

$wb = $excel.ActiveWorkbook

$wsMS = $wb.Worksheets.Item("MSData")
$wsFlux = $wb.Worksheets.Item("FluxData")
$wsTracer = $wb.Worksheets.Item("TracerData")

# --- FluxData: insert the new "EX_glc__D_e.f" flux row -------------------
# Row 3 currently holds "EX_c5sugal_e.f"; push it (and everything below)
# down one row, then populate the freshly-opened row 3.
[void]$wsFlux.Rows.Item(3).Insert()

$wsFlux.Range("A3").Value2 = "EX_glc__D_e.f"
$wsFlux.Range("B3").Value2 = 5.78096107699413
$wsFlux.Range("C3").Value2 = 3.01121981587178
$wsFlux.Rows.Item(3).RowHeight = 13.8

# The row that used to be row 3 ("EX_c5sugal_e.f") is now row 4; refresh
# its value/error.
$wsFlux.Range("B4").Value2 = 0.086805555555556
$wsFlux.Range("C4").Value2 = 0.148063493425702

# BIOMASS.f (row 2) error 1E-06 -> 0.0001
$wsFlux.Range("C2").Value2 = 0.0001

# Append a brand-new trailing row 38 ("DIL_ade_d1.f") mirroring the
# standard placeholder value/error pair used throughout the sheet.
$wsFlux.Range("A38").Value2 = "DIL_ade_d1.f"
$wsFlux.Range("B38").Value2 = 100
$wsFlux.Range("C38").Value2 = 0.0001

# --- View-state: zoom, active tab/sheet, selections -----------------------
[void]$wsMS.Activate()
$excel.ActiveWindow.Zoom = 95
[void]$wsMS.Range("A2").Select()

[void]$wsTracer.Activate()
$excel.ActiveWindow.Zoom = 95
[void]$wsTracer.Range("B1").Select()

# FluxData is activated last so it ends up the workbook's active tab/sheet.
[void]$wsFlux.Activate()
$excel.ActiveWindow.Zoom = 95
[void]$wsFlux.Range("A8").Select()
